# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Membrillo" (Mercado Mayorista Lo
# Valledor de Santiago) above the former row 122, pushing the existing
# rows 122-141 down to 124-143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 122 (this shifts old rows 122..141 down to 124..143)
$ws.Rows.Item(122).Insert()
$ws.Rows.Item(122).Insert()

# --- New row 122 ---
$ws.Cells.Item(122,1).Value  = 6
$ws.Cells.Item(122,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(122,3).Value  = "Metropolitana"
$ws.Cells.Item(122,4).Value  = 45015
$ws.Cells.Item(122,5).Value  = 13
$ws.Cells.Item(122,6).Value  = "Fruta"
$ws.Cells.Item(122,7).Value  = 100104
$ws.Cells.Item(122,8).Value  = "Frutos de pepita"
$ws.Cells.Item(122,9).Value  = 100104003
$ws.Cells.Item(122,10).Value = "Membrillo"
$ws.Cells.Item(122,11).Value = "Champion"
$ws.Cells.Item(122,12).Value = "Primera"
$ws.Cells.Item(122,13).Value = 8
$ws.Cells.Item(122,14).Value = 230000
$ws.Cells.Item(122,15).Value = 230000
$ws.Cells.Item(122,16).Value = 230000
$ws.Cells.Item(122,17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(122,18).Value = "Región de O'Higgins"
$ws.Cells.Item(122,19).Value = 511
$ws.Cells.Item(122,20).Value = 450

# --- New row 123 ---
$ws.Cells.Item(123,1).Value  = 6
$ws.Cells.Item(123,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(123,3).Value  = "Metropolitana"
$ws.Cells.Item(123,4).Value  = 45015
$ws.Cells.Item(123,5).Value  = 13
$ws.Cells.Item(123,6).Value  = "Fruta"
$ws.Cells.Item(123,7).Value  = 100104
$ws.Cells.Item(123,8).Value  = "Frutos de pepita"
$ws.Cells.Item(123,9).Value  = 100104003
$ws.Cells.Item(123,10).Value = "Membrillo"
$ws.Cells.Item(123,11).Value = "Champion"
$ws.Cells.Item(123,12).Value = "Segunda"
$ws.Cells.Item(123,13).Value = 10
$ws.Cells.Item(123,14).Value = 200000
$ws.Cells.Item(123,15).Value = 200000
$ws.Cells.Item(123,16).Value = 200000
$ws.Cells.Item(123,17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(123,18).Value = "Región de O'Higgins"
$ws.Cells.Item(123,19).Value = 444
$ws.Cells.Item(123,20).Value = 450
